$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the NB row (was row 8); SVM's row (was row 9) shifts up to become row 8 ---
$ws.Rows(8).Delete()

# --- The shifted-up row's index value (A8) still reads the old "7"; reset to 6 ---
$ws.Range("A8").Value = 6

# --- Extend header row: copy style from existing header cells into the new H1:L1 range ---
$ws.Range("C1:G1").Copy($ws.Range("H1"))

# --- Header text (row 1): split each "<period> Base" header into mean/std pairs ---
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# --- Algorithm name change (row 5 was CART, now DTREE) ---
$ws.Range("B5").Value = "DTREE"

# --- Row 2: LR ---
$ws.Range("C2").Value = 0.8290063035045472
$ws.Range("D2").Value = 0.01710508836750712
$ws.Range("E2").Value = 0.8084703364786201
$ws.Range("F2").Value = 0.01414886450690381
$ws.Range("G2").Value = 0.7895844436119595
$ws.Range("H2").Value = 0.02115744017527282
$ws.Range("I2").Value = 0.7727687108452517
$ws.Range("J2").Value = 0.03426160542212195
$ws.Range("K2").Value = 0.7615325412964358
$ws.Range("L2").Value = 0.02484367728360461

# --- Row 3: LDA ---
$ws.Range("C3").Value = 0.8328326018141052
$ws.Range("D3").Value = 0.01453236676767622
$ws.Range("E3").Value = 0.8204409948972897
$ws.Range("F3").Value = 0.01184206501536745
$ws.Range("G3").Value = 0.7989486317172515
$ws.Range("H3").Value = 0.02073440021578505
$ws.Range("I3").Value = 0.7870115079347577
$ws.Range("J3").Value = 0.03155701561389206
$ws.Range("K3").Value = 0.7676355213881721
$ws.Range("L3").Value = 0.02374883794081912

# --- Row 4: KNN ---
$ws.Range("C4").Value = 0.7687026036942916
$ws.Range("D4").Value = 0.01418735014209643
$ws.Range("E4").Value = 0.7586275181932673
$ws.Range("F4").Value = 0.02173391590076622
$ws.Range("G4").Value = 0.7692410224457522
$ws.Range("H4").Value = 0.02317959212935706
$ws.Range("I4").Value = 0.781363427111565
$ws.Range("J4").Value = 0.0211363178825121
$ws.Range("K4").Value = 0.7635984706755519
$ws.Range("L4").Value = 0.01246112832963159

# --- Row 5: DTREE ---
$ws.Range("C5").Value = 0.7640499797124396
$ws.Range("D5").Value = 0.03388293910715554
$ws.Range("E5").Value = 0.7543945795039878
$ws.Range("F5").Value = 0.02218091556870864
$ws.Range("G5").Value = 0.7598917605137633
$ws.Range("H5").Value = 0.02521632955632306
$ws.Range("I5").Value = 0.7411854895310889
$ws.Range("J5").Value = 0.02520780971262462
$ws.Range("K5").Value = 0.7328947042320806
$ws.Range("L5").Value = 0.01968859943852495

# --- Row 6: RTREE ---
$ws.Range("C6").Value = 0.7673980832291649
$ws.Range("D6").Value = 0.01949207555067016
$ws.Range("E6").Value = 0.7763326088178968
$ws.Range("F6").Value = 0.02474204593949417
$ws.Range("G6").Value = 0.7647255713539701
$ws.Range("H6").Value = 0.02360439616417246
$ws.Range("I6").Value = 0.7493663692847314
$ws.Range("J6").Value = 0.01999342872752569
$ws.Range("K6").Value = 0.7342399151372088
$ws.Range("L6").Value = 0.0249900055818145

# --- Row 7: XTREE ---
$ws.Range("C7").Value = 0.8312137226594478
$ws.Range("D7").Value = 0.01182619287020923
$ws.Range("E7").Value = 0.8197744413807424
$ws.Range("F7").Value = 0.01902258772529293
$ws.Range("G7").Value = 0.8039360545542816
$ws.Range("H7").Value = 0.02162671793513873
$ws.Range("I7").Value = 0.8059093083794024
$ws.Range("J7").Value = 0.02668171393622859
$ws.Range("K7").Value = 0.7827149034987756
$ws.Range("L7").Value = 0.02746805779067175

# --- Row 8: SVM ---
$ws.Range("C8").Value = 0.8337531560915117
$ws.Range("D8").Value = 0.01205233928829273
$ws.Range("E8").Value = 0.8226755030536056
$ws.Range("F8").Value = 0.01450328094857435
$ws.Range("G8").Value = 0.8180065468735366
$ws.Range("H8").Value = 0.02201438389461094
$ws.Range("I8").Value = 0.8079560318322775
$ws.Range("J8").Value = 0.02920561019688661
$ws.Range("K8").Value = 0.7961150438724658
$ws.Range("L8").Value = 0.02121312981701478

Write-Output "done"
